# chore: update Sheets via scheduled runner
#
# Refreshes the cached market-board derived figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ, columns H:N) on the per-job Leve profit sheets.
# These are plain scraped values (no formulas), so the refresh is just a set
# of direct cell writes; a couple of rows also gained/lost a cached profit
# figure entirely (ARM rows 64/67 lost their HQ profit figure, CUL row 40 and
# LTW rows 22/27 gained one) since that market data point went missing/found.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2937.914
$ws.Range("I138").Value = 1311.6052
$ws.Range("J138").Value = 4061.5454
$ws.Range("K138").Value = 3934.8156
$ws.Range("L138").Value = 12184.6362
$ws.Range("M138").Value = 1205.1844
$ws.Range("N138").Value = -22464.6362

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""
$ws.Range("H74").Value = 4875.243
$ws.Range("I74").Value = 3162.1738
$ws.Range("J74").Value = 7689.5713
$ws.Range("K74").Value = 3162.1738
$ws.Range("L74").Value = 7689.5713
$ws.Range("M74").Value = -2288.1738
$ws.Range("N74").Value = -9437.5713
$ws.Range("H77").Value = 4875.243
$ws.Range("I77").Value = 3162.1738
$ws.Range("J77").Value = 7689.5713
$ws.Range("K77").Value = 15810.869
$ws.Range("L77").Value = 38447.85649999999
$ws.Range("M77").Value = -11442.869
$ws.Range("N77").Value = -47183.85649999999
$ws.Range("H102").Value = 2616.182
$ws.Range("I102").Value = 2439.8572
$ws.Range("K102").Value = 2439.8572
$ws.Range("M102").Value = -817.8571999999999
$ws.Range("H132").Value = 5670.772
$ws.Range("I132").Value = 4599.276
$ws.Range("J132").Value = 6780.5356
$ws.Range("K132").Value = 13797.828
$ws.Range("L132").Value = 20341.6068
$ws.Range("M132").Value = -11267.828
$ws.Range("N132").Value = -25401.6068

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 32590.5
$ws.Range("J62").Value = 35181
$ws.Range("L62").Value = 35181
$ws.Range("N62").Value = -36553
$ws.Range("H65").Value = 32590.5
$ws.Range("J65").Value = 35181
$ws.Range("L65").Value = 105543
$ws.Range("N65").Value = -112407
$ws.Range("H105").Value = 4502.851
$ws.Range("I105").Value = 4074.9
$ws.Range("K105").Value = 4074.9
$ws.Range("M105").Value = -2327.9

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1044.7858
$ws.Range("I16").Value = 900.2222
$ws.Range("J16").Value = 1305
$ws.Range("K16").Value = 900.2222
$ws.Range("L16").Value = 1305
$ws.Range("M16").Value = -613.2222
$ws.Range("N16").Value = -1879
$ws.Range("H22").Value = 129.57143
$ws.Range("I22").Value = 106.75
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = 106.75
$ws.Range("L22").Value = 160
$ws.Range("M22").Value = 243.25
$ws.Range("N22").Value = -860
$ws.Range("H31").Value = 2031.1864
$ws.Range("I31").Value = 1659.2703
$ws.Range("J31").Value = 2656.682
$ws.Range("K31").Value = 1659.2703
$ws.Range("L31").Value = 2656.682
$ws.Range("M31").Value = -1364.2703
$ws.Range("N31").Value = -3246.682
$ws.Range("H34").Value = 2031.1864
$ws.Range("I34").Value = 1659.2703
$ws.Range("J34").Value = 2656.682
$ws.Range("K34").Value = 1659.2703
$ws.Range("L34").Value = 2656.682
$ws.Range("M34").Value = -1457.2703
$ws.Range("N34").Value = -3060.682
$ws.Range("H53").Value = 44800
$ws.Range("J53").Value = 44800
$ws.Range("L53").Value = 44800
$ws.Range("N53").Value = -46014
$ws.Range("H113").Value = 1044.7858
$ws.Range("I113").Value = 900.2222
$ws.Range("J113").Value = 1305
$ws.Range("K113").Value = 900.2222
$ws.Range("L113").Value = 1305
$ws.Range("M113").Value = 1269.7778
$ws.Range("N113").Value = -5645

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1550
$ws.Range("J22").Value = 1840
$ws.Range("L22").Value = 5520
$ws.Range("N22").Value = -5858
$ws.Range("H27").Value = 1550
$ws.Range("J27").Value = 1840
$ws.Range("L27").Value = 5520
$ws.Range("N27").Value = -5724
$ws.Range("H40").Value = 80
$ws.Range("I40").Value = 77.77778000000001
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 311.11112
$ws.Range("L40").Value = 400
$ws.Range("M40").Value = -242.11112
$ws.Range("N40").Value = -538
$ws.Range("H68").Value = 6506.9414
$ws.Range("I68").Value = 475.16666
$ws.Range("J68").Value = 9797
$ws.Range("K68").Value = 1425.49998
$ws.Range("L68").Value = 29391
$ws.Range("M68").Value = -614.4999800000001
$ws.Range("N68").Value = -31013
$ws.Range("H69").Value = 1514.4828
$ws.Range("J69").Value = 1631.5385
$ws.Range("L69").Value = 4894.6155
$ws.Range("N69").Value = -6516.6155
$ws.Range("H71").Value = 6506.9414
$ws.Range("I71").Value = 475.16666
$ws.Range("J71").Value = 9797
$ws.Range("K71").Value = 4276.49994
$ws.Range("L71").Value = 88173
$ws.Range("M71").Value = -220.4999399999997
$ws.Range("N71").Value = -96285
$ws.Range("H72").Value = 1514.4828
$ws.Range("J72").Value = 1631.5385
$ws.Range("L72").Value = 14683.8465
$ws.Range("N72").Value = -22795.8465
$ws.Range("H117").Value = 953.2222
$ws.Range("I117").Value = 333
$ws.Range("J117").Value = 1263.3334
$ws.Range("K117").Value = 999
$ws.Range("L117").Value = 3790.0002
$ws.Range("M117").Value = 2443
$ws.Range("N117").Value = -10674.0002
$ws.Range("H120").Value = 7110.524
$ws.Range("I120").Value = 5900
$ws.Range("J120").Value = 7395.353
$ws.Range("K120").Value = 17700
$ws.Range("L120").Value = 22186.059
$ws.Range("M120").Value = -12862
$ws.Range("N120").Value = -31862.059
$ws.Range("H123").Value = 7333.1333
$ws.Range("I123").Value = 7000
$ws.Range("J123").Value = 7356.9287
$ws.Range("K123").Value = 21000
$ws.Range("L123").Value = 22070.7861
$ws.Range("M123").Value = -18550
$ws.Range("N123").Value = -26970.7861
$ws.Range("H126").Value = 2631.7727
$ws.Range("I126").Value = 852.2222
$ws.Range("J126").Value = 3863.7693
$ws.Range("K126").Value = 2556.6666
$ws.Range("L126").Value = 11591.3079
$ws.Range("M126").Value = 2383.3334
$ws.Range("N126").Value = -21471.3079
$ws.Range("H129").Value = 1762.2693
$ws.Range("I129").Value = 1706.7778
$ws.Range("J129").Value = 1791.6471
$ws.Range("K129").Value = 5120.3334
$ws.Range("L129").Value = 5374.9413
$ws.Range("M129").Value = -120.3334000000004
$ws.Range("N129").Value = -15374.9413

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 957.1429000000001
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 957.1429000000001
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -1214
$ws.Range("H30").Value = 7504
$ws.Range("H35").Value = 25809.092
$ws.Range("I35").Value = 11515.5
$ws.Range("J35").Value = 28985.445
$ws.Range("K35").Value = 11515.5
$ws.Range("L35").Value = 28985.445
$ws.Range("M35").Value = -11179.5
$ws.Range("N35").Value = -29657.445
$ws.Range("H100").Value = 3651.2273
$ws.Range("I100").Value = 2252.25
$ws.Range("J100").Value = 5330
$ws.Range("K100").Value = 2252.25
$ws.Range("L100").Value = 5330
$ws.Range("M100").Value = -1711.25
$ws.Range("N100").Value = -6412
